$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "weeknr 49web" sheet by copying the "weeknr 48web" sheet
#    (same column layout / template as the other "...web" weeks) and placing
#    it right after "weeknr 49game", just before "totaal".
# ---------------------------------------------------------------------------
$src   = $wb.Worksheets.Item("weeknr 48web")
$after = $wb.Worksheets.Item("weeknr 49game")
$src.Copy($null, $after)

$newSheet = $wb.Worksheets.Item("weeknr 48web (2)")
$newSheet.Name = "weeknr 49web"

# ---------------------------------------------------------------------------
# 2) Update the "totaal" sheet: label the new row ("49 web") before we touch
#    any other new text, so the shared-string table appends "49 web" first.
# ---------------------------------------------------------------------------
$totaal = $wb.Worksheets.Item("totaal")
$totaal.Range("A9").Value2 = "49 web"

# ---------------------------------------------------------------------------
# 3) Fill in the first logged activity on the new "weeknr 49web" sheet.
# ---------------------------------------------------------------------------
$newSheet.Range("A8").Value2 = "Maandag"
$newSheet.Range("B8").Value2 = 41610
$newSheet.Range("C8").Value2 = 0.61319444444444449
$newSheet.Range("D8").Value2 = 0.61805555555555558
$newSheet.Range("F8").Value2 = "aptana, wamp opstarten"

# Rows 9:18 are still an empty template (only row/number placeholders) -
# clear the copied sample times & activity text but keep their formatting.
$newSheet.Range("C9:D18").ClearContents()
$newSheet.Range("F9:F18").ClearContents()

# ---------------------------------------------------------------------------
# 4) Point the "totaal" sheet's new row at the new sheet's grand total.
# ---------------------------------------------------------------------------
$totaal.Range("B9").Formula = "='weeknr 49web'!H31"

# ---------------------------------------------------------------------------
# 5) "weeknr 49game": the Monday entry actually happened a week later.
# ---------------------------------------------------------------------------
$game49 = $wb.Worksheets.Item("weeknr 49game")
$game49.Range("B8").Value2 = 41610

# ---------------------------------------------------------------------------
# 6) Selections / active sheet, matching the saved view state.
# ---------------------------------------------------------------------------
$game49.Range("B9").Select()
$totaal.Range("B10").Select()

$newSheet.Activate()
$newSheet.Range("C9").Select()
